# edit.ps1
# Applies the "add uuid, fix delete function" change to ribao.xlsx
#
# For every sheet a new "UUID" column is inserted right after the date
# column (i.e. becomes column C, pushing the former C.. columns one to
# the right). Several sheets also have their (now stale / sample) data
# rows removed, keeping only the header row. The "other / 其它" sheet
# keeps its sample row but refreshes it with a new date, a UUID value and
# an extra content column. Various cosmetic view/outline bits are also
# refreshed to match.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to be stored as literal text (avoids Excel turning
    # all-digit / date-shaped strings into numbers or dates), while not
    # leaving a lingering "quote prefix" number format behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 1: 渗透测试
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns("C:C").Insert()
$ws1.Range("C1").Value = "UUID"

Set-TextValue $ws1.Range("A2") "ZhangShan"
Set-TextValue $ws1.Range("B2") "2024-04-01"
Set-TextValue $ws1.Range("C2") "e180c0d7-1485-4200-ab8e-dbb088a2bb97"
Set-TextValue $ws1.Range("D2") "系统名称1"
Set-TextValue $ws1.Range("E2") "版本线"
Set-TextValue $ws1.Range("F2") "11111111"
Set-TextValue $ws1.Range("G2") "222222222"

$ws1.Outline.ShowLevels(0, 0)
$ws1.Range("D5").Select()

# ---------------------------------------------------------------------
# Sheet 2: 渗透复测
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Columns("C:C").Insert()
$ws2.Range("C1").Value = "UUID"

$ws2.Outline.ShowLevels(0, 0)
$ws2.Range("D7").Select()

# ---------------------------------------------------------------------
# Sheet 3: 代码审计
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Columns("C:C").Insert()
$ws3.Range("C1").Value = "UUID"
$ws3.Rows("2:3").Delete()

$ws3.Outline.ShowLevels(0, 0)
$ws3.Range("F13").Select()

# ---------------------------------------------------------------------
# Sheet 4: 漏洞审核
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Columns("C:C").Insert()
$ws4.Range("C1").Value = "UUID"
$ws4.Rows("2:2").Delete()

$ws4.Outline.ShowLevels(0, 0)
$ws4.Range("H20").Select()

# ---------------------------------------------------------------------
# Sheet 5: 安全开发
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Columns("C:C").Insert()
$ws5.Range("C1").Value = "UUID"
$ws5.Rows("2:2").Delete()

$ws5.Outline.ShowLevels(0, 0)
$ws5.Range("C5").Select()

# ---------------------------------------------------------------------
# Sheet 6: 文档编制
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Columns("C:C").Insert()
$ws6.Range("C1").Value = "UUID"
$ws6.Rows("2:3").Delete()

$ws6.Outline.ShowLevels(0, 0)
$ws6.Range("E9").Select()

# ---------------------------------------------------------------------
# Sheet 7: 其它
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)
$ws7.Columns("C:C").Insert()
$ws7.Range("C1").Value = "UUID"

Set-TextValue $ws7.Range("B2") "2024-04-01"
Set-TextValue $ws7.Range("C2") "b1f13017-a213-4a9f-9dbe-c40644dd8f7a"
Set-TextValue $ws7.Range("D2") "2222222222"

# This sheet is the one left active/selected in the saved workbook, so
# select it last.
$ws7.Range("H11").Select()
